# Change Slide 1's layout from "Title and Content" to "Title Slide".
# This re-types the title placeholder to ctrTitle and the content
# placeholder to subTitle (idx=1), while keeping the existing text.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$d = $p.Designs.Item(1)
$titleSlideLayout = $d.SlideMaster.CustomLayouts.Item(1)

# Remember the current text of the two placeholders before switching
# layouts (the old shapes stay put, unfilled new ones get added).
$oldTitleShape = $s.Shapes.Item(1)
$oldBodyShape = $s.Shapes.Item(2)
$titleText = $oldTitleShape.TextFrame.TextRange.Text
$bodyText = $oldBodyShape.TextFrame.TextRange.Text

$s.CustomLayout = $titleSlideLayout

# After the assignment the slide has 4 shapes: the two original
# (still typed title / idx=1, still holding the old text) followed by
# two freshly-added empty placeholders from the new layout (ctrTitle,
# subTitle idx=1). Move the text across, then remove the stale shapes.
$newTitleShape = $s.Shapes.Item(3)
$newBodyShape = $s.Shapes.Item(4)
$newTitleShape.TextFrame.TextRange.Text = $titleText
$newBodyShape.TextFrame.TextRange.Text = $bodyText

# Deleting a stale placeholder first resets it to an empty layout
# placeholder (a "ghost"); deleting that ghost removes it for good.
$s.Shapes.Item(1).Delete()
$s.Shapes.Item(1).Delete()
$s.Shapes.Item(1).Delete()
$s.Shapes.Item(1).Delete()
